$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 1.349324843968644
$ws.Range("AG2").Value = -6.110329073966087
$ws.Range("AD3").Value = 0.1790125063937677
$ws.Range("AG3").Value = -12.55213029169902
$ws.Range("AD4").Value = 0.3251491368039667
$ws.Range("AG4").Value = -6.421233489441632
$ws.Range("AD5").Value = 1.295364897130082
$ws.Range("AG5").Value = -5.262177032631001
$ws.Range("AD6").Value = 0.2573768658608793
$ws.Range("AG6").Value = -11.37043260665906
$ws.Range("AD7").Value = 0.6643901723454653
$ws.Range("AG7").Value = -8.738976210705758
$ws.Range("AD8").Value = 0.0154326596133041
$ws.Range("AG8").Value = -13.29114803329192
$ws.Range("AD9").Value = 0.5184410162249959
$ws.Range("AG9").Value = -8.72907618319705
$ws.Range("AD10").Value = 0.3700482611681618
$ws.Range("AG10").Value = -10.59589935263373
$ws.Range("AD11").Value = 0.2454642123782856
$ws.Range("AG11").Value = -10.47658916606095
$ws.Range("AD12").Value = 0.05169514222383378
$ws.Range("AG12").Value = -10.23999818360474
$ws.Range("AD13").Value = 0.3448120709288106
$ws.Range("AG13").Value = -10.07268161101236
$ws.Range("AD14").Value = -0.2355492703904264
$ws.Range("AG14").Value = -6.820554729806321
$ws.Range("AD15").Value = 0.9472263598682301
$ws.Range("AG15").Value = -8.139702232314026
$ws.Range("AD16").Value = 0.2141670668851253
$ws.Range("AG16").Value = -9.14915920128017
$ws.Range("AD17").Value = 0.5258243701171119
$ws.Range("AG17").Value = -12.05258875414952
$ws.Range("AD18").Value = 0.1833057836558002
$ws.Range("AG18").Value = -10.95016346752053
$ws.Range("AD19").Value = 0.1283109762246061
$ws.Range("AG19").Value = -10.9523170000961
$ws.Range("AD20").Value = 0.2181686601032652
$ws.Range("AG20").Value = -11.14750700394278
$ws.Range("AD21").Value = 0.12262271326239
$ws.Range("AG21").Value = -10.5956270662656
$ws.Range("AD22").Value = 0.08755629789200478
$ws.Range("AG22").Value = -8.869760543648585
$ws.Range("AD23").Value = 0.2539916717785259
$ws.Range("AG23").Value = -8.196066117660484
$ws.Range("AD24").Value = 0.9667243371685512
$ws.Range("AG24").Value = -8.019025729421545
$ws.Range("AD25").Value = 0.8735024107446513
$ws.Range("AG25").Value = -5.278420909372512
$ws.Range("AD26").Value = 0.6485321294332147
$ws.Range("AG26").Value = -11.96251076101962
$ws.Range("AD27").Value = 0.7147900127387261
$ws.Range("AG27").Value = -9.273139604613139
$ws.Range("AD28").Value = 1.21651060392319
$ws.Range("AG28").Value = -11.36309339723514
$ws.Range("AD29").Value = -0.1278783053910493
$ws.Range("AG29").Value = -9.815400723727365
$ws.Range("AD30").Value = -1.28158685972077
$ws.Range("AG30").Value = -10.54341820465791
$ws.Range("AD31").Value = 0.7622049139410226
$ws.Range("AG31").Value = -9.279090242039002
$ws.Range("AD32").Value = -0.7800801661261121
$ws.Range("AG32").Value = -10.09378773240133
$ws.Range("AD33").Value = 0.6516771192456889
$ws.Range("AG33").Value = -8.959699695999088
$ws.Range("AD34").Value = 0.2036337351286908
$ws.Range("AG34").Value = -8.647342028124038
$ws.Range("AD35").Value = 0.7926442020663729
$ws.Range("AG35").Value = -9.943078560421824
$ws.Range("AD36").Value = 0.6496091768848999
$ws.Range("AG36").Value = -11.93894908408843
$ws.Range("AD37").Value = 0.05582418165583897
$ws.Range("AG37").Value = -9.45192088529968
$ws.Range("AD38").Value = 0.2055341294371552
$ws.Range("AG38").Value = -9.976478896918959
$ws.Range("AD39").Value = 0.04042874502368653
$ws.Range("AG39").Value = -12.73155676084157
$ws.Range("AD40").Value = 1.581307269299106
$ws.Range("AG40").Value = -6.268205720733754
$ws.Range("AD41").Value = 0.05606475215748565
$ws.Range("AG41").Value = -12.84814564646513
$ws.Range("AD42").Value = 0.5120571413802395
$ws.Range("AG42").Value = -10.87016467859482
$ws.Range("AD43").Value = 1.171238522772072
$ws.Range("AG43").Value = -8.925435614314068
$ws.Range("AD44").Value = 0.04952614163671054
$ws.Range("AG44").Value = -10.42783355760777
$ws.Range("AD45").Value = 0.7465755421807339
$ws.Range("AG45").Value = -8.019222898608213
$ws.Range("AD46").Value = 0.6797772062528168
$ws.Range("AG46").Value = -10.94557452530746
$ws.Range("AD47").Value = 0.3458254045777432
$ws.Range("AG47").Value = -10.84139785043272
$ws.Range("AD48").Value = 0.02688030345492576
$ws.Range("AG48").Value = -11.18451515084376
$ws.Range("AD49").Value = 0.04137663367350292
$ws.Range("AG49").Value = -8.677714952258562
$ws.Range("AD50").Value = 0.6863028741057466
$ws.Range("AG50").Value = -7.413685313075046
$ws.Range("AD51").Value = 0.1862437341284388
$ws.Range("AG51").Value = -8.230410998748049
$ws.Range("AD52").Value = 0.4312064758209718
$ws.Range("AG52").Value = -10.47379969847979
$ws.Range("AD53").Value = -0.1669175544610551
$ws.Range("AG53").Value = -7.817074615968632
$ws.Range("AD54").Value = 0.04286757755881693
$ws.Range("AG54").Value = -8.432808342726089
$ws.Range("AD55").Value = 0.06375664783978777
$ws.Range("AG55").Value = -9.612153771096022
$ws.Range("AD56").Value = -0.02842617464579093
$ws.Range("AG56").Value = -10.10959910489474
$ws.Range("AD57").Value = 0.1535286260219229
$ws.Range("AG57").Value = -9.190935359354258
$ws.Range("AD58").Value = -0.002499888651234683
$ws.Range("AG58").Value = -7.860314983029536
$ws.Range("AD59").Value = -0.651592177382373
$ws.Range("AG59").Value = -11.5130204287197
$ws.Range("AD60").Value = 0.1918172034293323
$ws.Range("AG60").Value = -11.47415877205243
